$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Add the K/24, MOD(K,24), MOD(L,24), MOD(M,24), MOD(N,24) helper
#    formulas in columns Q:U for rows 2 through 17 (departure/arrival
#    hour-of-day helper columns).
# ------------------------------------------------------------------
for ($r = 2; $r -le 17; $r++) {
    $ws.Range("Q$r").Formula = "=K$r/24"
    $ws.Range("R$r").Formula = "=MOD(K$r,24)"
    $ws.Range("S$r").Formula = "=MOD(L$r,24)"
    $ws.Range("T$r").Formula = "=MOD(M$r,24)"
    $ws.Range("U$r").Formula = "=MOD(N$r,24)"
}

# ------------------------------------------------------------------
# 2) Append a new vehicle row (Truck35, route Neuss -> Nuremberg),
#    following the same pattern as the other Truck rows above it.
#    Columns E/F/G store the speed/cost figures as TEXT (matching the
#    rest of the table), so the numeric-looking values are routed
#    through a scratch cell + paste-values so Excel keeps them as
#    shared-string text instead of silently coercing them to numbers.
# ------------------------------------------------------------------
$ws.Range("A118").Value = "Truck35"
$ws.Range("B118").Value = 1000
$ws.Range("C118").Value = 1000
$ws.Range("D118").Value = 75

$ws.Range("Z1").Formula = '="30.98"'
$ws.Range("Z1").Copy()
$ws.Range("E118").PasteSpecial(-4163)

$ws.Range("Z1").Formula = '="0.2758"'
$ws.Range("Z1").Copy()
$ws.Range("F118").PasteSpecial(-4163)

$ws.Range("Z1").Formula = '="0.8866"'
$ws.Range("Z1").Copy()
$ws.Range("G118").PasteSpecial(-4163)

$ws.Range("Z1").ClearContents()

$ws.Range("H118").Value = 3
$ws.Range("I118").Value = "Neuss"
$ws.Range("J118").Value = "Nuremberg"
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 10000
$ws.Range("M118").Value = 0
$ws.Range("N118").Value = 10000

# ------------------------------------------------------------------
# 3) Refresh the view: zoom level + final selection/scroll position,
#    matching where the user ended up after adding the new row.
# ------------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 85
$ws.Application.ActiveWindow.ScrollRow = 105
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("I128").Select() | Out-Null
